# Update speaker/content text on Slide 1 (Redux overview) and Slide 2
# (Redux best-practices / Context API comparison) content placeholders.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 - "Content Placeholder 2"
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$shape1 = $slide1.Shapes.Item(2)
$tr1 = $shape1.TextFrame.TextRange

$tr1.Paragraphs(2, 1).Runs(1, 1).Text = "What is Redux? - A state management library, centralizes application state, works with React but is library-agnostic."
$tr1.Paragraphs(3, 1).Runs(1, 1).Text = "Why Use Redux? - Manages complex state logic, centralizes data, improves predictability and debugging."
$tr1.Paragraphs(6, 1).Runs(1, 1).Text = "Redux Flow - User interaction → Dispatch action → Reducer updates state → UI re-renders."
$tr1.Paragraphs(8, 1).Runs(1, 1).Text = "Modern Redux (RTK) - Simplifies Redux logic, reduces boilerplate, includes createSlice and createAsyncThunk."
$tr1.Paragraphs(9, 1).Runs(1, 1).Text = "When NOT to Use Redux? - Small apps, simple state, or when Context API/useState suffices."
$tr1.Paragraphs(10, 1).Runs(1, 1).Text = "Benefits - Centralized state, predictable transitions, easy debugging, scalable architecture."
$tr1.Paragraphs(12, 1).Runs(1, 1).Text = "Best Practices - Pure reducers, normalized state, avoid non-serializable data."

# ---------------------------------------------------------------------
# Slide 2 - "Content Placeholder 2"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange

$tr2.Paragraphs(2, 1).Runs(1, 1).Text = "Use Redux Toolkit to simplify setup."
$tr2.Paragraphs(3, 1).Runs(1, 1).Text = "Use selectors to avoid tight coupling to state shape."
$tr2.Paragraphs(4, 1).Runs(1, 1).Text = "Comparison of Redux and Context API:"
$tr2.Paragraphs(5, 1).Runs(1, 1).Text = "- Redux: Best for large apps, complex state flows."
$tr2.Paragraphs(6, 1).Runs(1, 1).Text = "- Context API: Best for small to medium apps, simple state."

$newParas = "`rWhen to use Context API:" + `
            "`r- For static or lightweight state (e.g., theme, auth)." + `
            "`r- When avoiding extra libraries." + `
            "`r- When state changes infrequently." + `
            "`rWhen to use Redux:" + `
            "`r- For large, shared, or complex app state." + `
            "`r- When middleware, async handling, or devtools are needed." + `
            "`r- For predictable and testable state transitions."

$null = $tr2.InsertAfter($newParas)
